# Add season-record columns (Wins / Losses / Ties) to the roster table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells AD1:AF1 need the same formatting as the other header cells
# (bold, centered, bordered -- the style already used by A1:AC1). Copy that
# formatting across first, then stamp in the header text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-45) gets the same 1999 Baltimore Orioles season
# record: 78 wins, 84 losses, 0 ties.
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value = 78
    $ws.Cells.Item($r, 31).Value = 84
    $ws.Cells.Item($r, 32).Value = 0
}
